$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D so the existing "Tipo" column shifts to E
$ws.Range("D1").EntireColumn.Insert()

# New header for the inserted column - copy the header formatting from C1
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("D1").Value = "MAE"

# Update the numeric values in row 2
$ws.Range("B2").Value = 0.3602986410487848
$ws.Range("C2").Value = 0.9950205752246711
$ws.Range("D2").Value = 0.4980033758824289
